$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new "2020" column (Q) mirroring the existing "2019"
# column (P): same per-row formatting, new data values for each line.
# Copy column P's formatting into column Q first (keeps every row's
# number format / font / borders identical to its P neighbour), then
# overwrite the values row by row.
$ws.Range("P3:P37").Copy()
$ws.Range("Q3:Q37").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("Q3").Value = 2020

$ws.Range("Q4").Value = 0.1
$ws.Range("Q5").Value = 0.1
$ws.Range("Q6").Value = 0.1

$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("Q9").Value = 0

$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("Q12").Value = 0

$ws.Range("Q13").Value = 0
$ws.Range("Q14").Value = 0.1
$ws.Range("Q15").Value = 0

$ws.Range("Q16").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("Q18").Value = 0

$ws.Range("Q19").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("Q21").Value = 0

$ws.Range("Q22").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("Q24").Value = 0

$ws.Range("Q25").Value = 0.1
$ws.Range("Q26").Value = 0.2
$ws.Range("Q27").Value = 0.1

$ws.Range("Q28").Value = 0.3
$ws.Range("Q29").Value = 0.4
$ws.Range("Q30").Value = 0.2

$ws.Range("Q31").Value = 0.2
$ws.Range("Q32").Value = 0.2
$ws.Range("Q33").Value = 0.1

# Row 34 is a spacer row - keep it formatted but value-less, like P34.

$ws.Range("Q35").Value = 0
$ws.Range("Q36").Value = 0.1
$ws.Range("Q37").Value = 0.2

# Restore the cursor/selection that was active when the file was saved.
$ws.Range("P30").Select()
